# Generate Report for Handback
#
# The handback transform failed for the 5b3e9f2a-... file because the
# handback file name didn't match the handoff file name. Update the
# "Status" cells for that row on all three sheets, and record the error
# detail on the per-locale (zh-cn / de-de) sheets.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns for the 5b3e9f2a-... row
$ws_overview.Range("B7").Value = "Handback transform failed"
$ws_overview.Range("C7").Value = "Handback transform failed"

# zh-cn sheet: Status column + new Error Detail
$ws_zhcn.Range("C7").Value = "Handback transform failed"
$ws_zhcn.Range("L7").Value = "Handback file name: exgd1jxa.5sg is different with handoff file name: 5b3e9f2a-6635-4976-ae56-befec5dc6a7e.c42a801fe0e187ad4b0a7f0e3e0f27dfebce6199.zh-cn."

# de-de sheet: Status column + new Error Detail
$ws_dede.Range("C7").Value = "Handback transform failed"
$ws_dede.Range("L7").Value = "Handback file name: exgd1jxa.5sg is different with handoff file name: 5b3e9f2a-6635-4976-ae56-befec5dc6a7e.c42a801fe0e187ad4b0a7f0e3e0f27dfebce6199.de-de."
